# Apply the latest cryptos data refresh (prices & 1h volume changes) per GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Addr, $Val)
    $r = $ws.Range($Addr)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

Set-CellText "D2" "69.440.29"
Set-CellText "E2" "  +0.45%  "
Set-CellText "D3" "3.893.66"
Set-CellText "E3" "  +0.69%  "
Set-CellText "E4" "  +0.16%  "
Set-CellText "D5" "604.39"
Set-CellText "E5" "  +0.23%  "
Set-CellText "D6" "171.61"
Set-CellText "E6" "  +3.14%  "
Set-CellText "D7" "3.893.52"
Set-CellText "E7" "  +0.55%  "
Set-CellText "E8" "  +0.12%  "
Set-CellText "D9" "0.535"
Set-CellText "E9" "  +0.54%  "
Set-CellText "D10" "0.169"
Set-CellText "E10" "  +0.68%  "
Set-CellText "D11" "6.40"
Set-CellText "E11" "  +0.25%  "
Set-CellText "D12" "0.469"
Set-CellText "E12" "  +1.42%  "
Set-CellText "B13" "ShibaInu"
Set-CellText "C13" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-CellText "D13" "0.0000256"
Set-CellText "E13" "  +4.42%  "
Set-CellText "B14" "Avalanche"
Set-CellText "C14" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-CellText "D14" "38.35"
Set-CellText "E14" "  +3.11%  "
Set-CellText "D15" "4.542.91"
Set-CellText "E15" "  +0.92%  "
Set-CellText "D16" "3.896.94"
Set-CellText "E16" "  +1.27%  "
Set-CellText "D17" "69.471.05"
Set-CellText "E17" "  +0.58%  "
Set-CellText "D18" "18.76"
Set-CellText "E18" "  +9.11%  "
Set-CellText "D19" "7.66"
Set-CellText "E19" "  +1.31%  "
Set-CellText "E20" "  -0.63%  "
Set-CellText "D21" "10.98"
Set-CellText "E21" "  -2.72%  "
Set-CellText "D22" "488.38"
Set-CellText "E22" "  -0.21%  "
Set-CellText "D23" "0.749"
Set-CellText "E23" "  +3.46%  "
Set-CellText "D24" "0.0000167"
Set-CellText "E24" "  +1.18%  "
Set-CellText "D25" "85.34"
Set-CellText "E25" "  +1.28%  "
Set-CellText "D26" "2.30"
Set-CellText "E26" "  +1.50%  "
Set-CellText "D27" "12.46"
Set-CellText "E27" "  +2.46%  "
Set-CellText "D28" "10.14"
Set-CellText "E28" "  +0.72%  "
Set-CellText "E29" "  +0.03%  "
Set-CellText "D30" "3.00"
Set-CellText "E30" "  +1.69%  "
Set-CellText "D31" "4.038.00"
Set-CellText "D32" "7.88"
Set-CellText "E32" "  -0.87%  "
Set-CellText "D33" "2.36"
Set-CellText "E33" "  -0.92%  "
Set-CellText "D34" "31.81"
Set-CellText "E34" "  -0.90%  "
Set-CellText "D35" "3.857.00"
Set-CellText "E35" "  +1.45%  "
Set-CellText "D36" "0.107"
Set-CellText "E36" "  -0.59%  "
Set-CellText "B37" "dogwifhat"
Set-CellText "C37" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-CellText "D37" "3.41"
Set-CellText "E37" "  +11.31%  "
Set-CellText "B38" "Filecoin"
Set-CellText "C38" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-CellText "D38" "6.11"
Set-CellText "E38" "  +3.21%  "
Set-CellText "B39" "Kaspa"
Set-CellText "C39" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-CellText "D39" "0.142"
Set-CellText "E39" "  +0.95%  "
Set-CellText "E40" "  -0.94%  "
Set-CellText "E41" "  +0.18%  "
Set-CellText "D42" "0.324"
Set-CellText "E42" "  +1.43%  "
Set-CellText "D43" "2.09"
Set-CellText "E43" "  +4.41%  "
Set-CellText "D44" "434.37"
Set-CellText "E44" "  -2.71%  "
Set-CellText "E45" "  -1.21%  "
Set-CellText "D46" "8.69"
Set-CellText "E46" "  +2.46%  "
Set-CellText "E47" "  +0.01%  "
Set-CellText "D48" "0.000274"
Set-CellText "E48" "  +19.18%  "
Set-CellText "D49" "0.0363"
Set-CellText "E49" "  +1.14%  "
Set-CellText "D50" "40.01"
Set-CellText "E50" "  +1.66%  "
Set-CellText "D51" "141.60"
Set-CellText "E51" "  -0.36%  "
